$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Buy"
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
